$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.519.38'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.88%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.725.70'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.85%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.40'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.81%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("E8").Value = '  +1.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.753.06'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.113'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.08%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.391'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.13%  '

$ws.Range("E13").Value = '  +3.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.209.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.36'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.428.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.59%  '

$ws.Range("E17").Value = '  +7.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.745.44'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.63%  '

$ws.Range("E20").Value = '  +3.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '360.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.25%  '

$ws.Range("E22").Value = '  +1.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.538'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.78%  '

$ws.Range("E24").Value = '  -0.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.83%  '

$ws.Range("E26").Value = '  +4.13%  '

$ws.Range("E27").Value = '  +4.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.995'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0891'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.03'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.28%  '

$ws.Range("E31").Value = '  +7.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '174.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.46%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.22'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +18.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '20.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.84%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.77'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.45'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +9.64%  '

$ws.Range("E39").Value = '  +15.06%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '344.21'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.47%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.46%  '

$ws.Range("E43").Value = '  +7.83%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.74'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.91'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '140.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.52%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0588'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.32%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.645'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0256'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.80%  '

$ws.Range("E50").Value = '  +1.23%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.996'
$ws.Range("D51").Style = "Normal"

